$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 (shifts existing rows 7+ down by one)
$ws.Rows.Item(7).Insert()

# Populate the new row with the Address / adr bean entry
$ws.Range("B7").Value = "Address"
$ws.Range("C7").Value = "adr"

# Reflect the new selection state captured in the saved workbook
$ws.Range("C8").Select()
